# Updated Mail List 21/03/23. Maximising all the windows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old source/destination report file names (B4, B5)
$ws.Range("B4").ClearContents()
$ws.Range("B5").ClearContents()

# Update the folder path (B7) and turn it into a hyperlink
$ws.Range("B7").Value = "\\10.222.140.144\d\d\MIS\GLS\RPA_BOT\"
$ws.Hyperlinks.Add($ws.Range("B7"), "\\10.222.140.144\d\d\MIS\GLS\RPA_BOT\")

# Clear the statement date value (B8) but keep its style
$ws.Range("B8").ClearContents()

# Move the active selection
$ws.Range("B16").Select() | Out-Null

# Maximise the Excel application window
$excel.WindowState = -4137 | Out-Null
